$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'30.070.75"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.00%  "

$ws.Range("D3").Value = "'2.120.48"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.67%  "

$ws.Range("D4").Value = "'1.006"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.15%  "

$ws.Range("D5").Value = "'346.60"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.27%  "

$ws.Range("D6").Value = "'1.006"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.03%  "

$ws.Range("D7").Value = "'0.5199"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.35%  "

$ws.Range("D8").Value = "'0.4465"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.66%  "

$ws.Range("D9").Value = "'53.76"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +2.10%  "

$ws.Range("D10").Value = "'0.09384"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.44%  "

$ws.Range("D11").Value = "'1.184"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.57%  "

$ws.Range("D12").Value = "'25.39"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.18%  "

$ws.Range("D13").Value = "'8.576"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +5.34%  "

$ws.Range("B14").Value = "Polkadot"
$ws.Range("C14").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D14").Value = "'6.953"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +3.09%  "

$ws.Range("B15").Value = "WrappedEther"
$ws.Range("C15").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D15").Value = "'2.099.67"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.50%  "

$ws.Range("D16").Value = "'102.95"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +3.03%  "

$ws.Range("E17").Value = "  -0.48%  "

$ws.Range("E18").Value = "  -0.06%  "

$ws.Range("D19").Value = "'21.56"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +3.83%  "

$ws.Range("D20").Value = "'0.06688"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.27%  "

$ws.Range("E21").Value = "  +1.03%  "

$ws.Range("D22").Value = "'1.006"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.10%  "

$ws.Range("D23").Value = "'30.073.23"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.32%  "

$ws.Range("E24").Value = "  +0.30%  "

$ws.Range("D25").Value = "'2.317"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.95%  "

$ws.Range("D26").Value = "'22.12"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.10%  "

$ws.Range("D27").Value = "'2.543"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.39%  "

$ws.Range("D28").Value = "'162.64"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.06%  "

$ws.Range("D29").Value = "'134.09"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.17%  "

$ws.Range("D30").Value = "'1.155"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.18%  "

$ws.Range("D31").Value = "'1.794"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +9.54%  "

$ws.Range("D32").Value = "'0.1057"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.40%  "

$ws.Range("D33").Value = "'6.270"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.04%  "

$ws.Range("D34").Value = "'6.615"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +6.20%  "

$ws.Range("D35").Value = "'3.967"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.27%  "

$ws.Range("E36").Value = "  +6.09%  "

$ws.Range("D37").Value = "'0.02616"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.81%  "

$ws.Range("D38").Value = "'0.06869"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.02%  "

$ws.Range("D39").Value = "'0.7086"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.34%  "

$ws.Range("D40").Value = "'12.73"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.33%  "

$ws.Range("D41").Value = "'1.332"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.36%  "

$ws.Range("D42").Value = "'0.2246"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.91%  "

$ws.Range("D43").Value = "'0.6853"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.98%  "

$ws.Range("D44").Value = "'14.61"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +2.56%  "

$ws.Range("D45").Value = "'2.373"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +3.61%  "

$ws.Range("D46").Value = "'1.005"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.09%  "

$ws.Range("D47").Value = "'1.268"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +8.96%  "

$ws.Range("B48").Value = "PancakeSwap"
$ws.Range("C48").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D48").Value = "'3.629"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.31%  "

$ws.Range("B49").Value = "BabyDogeCoin"
$ws.Range("C49").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D49").Value = "'0.00000000356"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.74%  "

$ws.Range("D50").Value = "'1.227"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.22%  "

$ws.Range("D51").Value = "'83.40"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.80%  "
